$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.459.58"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").Value = "3.706.62"
$ws.Range("E3").Value = "  +8.39%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.45"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").Value = "3.697.60"
$ws.Range("E7").Value = "  +8.35%  "
$ws.Range("E8").Value = "  +4.29%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.609"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.22"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "4.303.79"
$ws.Range("E14").Value = "  +8.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "681.06"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").Value = "3.717.35"
$ws.Range("E17").Value = "  +8.47%  "
$ws.Range("D18").Value = "71.608.89"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.122"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.97"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.49"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +19.49%  "
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.36"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.22"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +11.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "591.66"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.20"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("E35").Value = "  +4.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.10"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "3.675.88"
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("D40").Value = "0.0₃0771"
$ws.Range("E40").Value = "  +6.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.32"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("E44").Value = "  +10.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.348"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.16%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.69%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.95"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.09%  "
